# "template master and sender master complete, general setting in progress"
#
# Target sheet is Sheet1 (the sheet that is tabSelected / ActiveSheet),
# a student/SMS "demo sample" roster with columns:
#   A: S.No.  B: RF ID Card No.  C: Student Name  D: Admission No.
#   E: Class & Section  F: Date Of Birth  G: Blood Group  ...

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column B ("RF ID Card No.") -----------------------------------------
# Rows 2 & 4 used to hold plain numbers (3 / 4); row 3 already held the
# text "5". Re-key all three RFID values as text "6", "7", "8" so the
# whole column is consistently text-formatted.
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "6"

$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "7"

$ws.Range("B4").NumberFormat = "@"
$ws.Range("B4").Value = "8"

# --- Column F ("Date Of Birth") -------------------------------------------
# Switch from the generic short-date format to an explicit dd/mm/yyyy
# display format across the header and all data rows, and update the
# second row's birth date to a new value (13-Dec-2010).
$ws.Range("F1:F4").NumberFormat = "dd\/mm\/yyyy"
$ws.Range("F2").Value = 40525

# --- Selection / view state -------------------------------------------
$ws.Range("F11").Select() | Out-Null

# --- Page setup (printing) ------------------------------------------------
$ws.PageSetup.PaperSize = 9   # xlPaperA4
$ws.PageSetup.Orientation = 1 # xlPortrait
